# Auto-generated edit script: updates crypto price/volume data per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.774.41"
$ws.Range("E2").Value = "  +3.09%  "

$ws.Range("D3").Value = "3.136.13"
$ws.Range("E3").Value = "  +2.18%  "

$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.02%  "

$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("D8").Value = "3.129.11"
$ws.Range("E8").Value = "  +2.40%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.537"
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.163"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +18.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.71"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.28%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.468"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.66%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000255"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.14%  "

$ws.Range("E15").Value = "  +0.15%  "

$ws.Range("D16").Value = "3.647.86"
$ws.Range("E16").Value = "  +2.07%  "

$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.71%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "63.696.93"
$ws.Range("E18").Value = "  +3.07%  "

$ws.Range("D19").Value = "3.130.13"
$ws.Range("E19").Value = "  +2.37%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "466.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.71%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.50%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.734"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.23%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.83%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("E26").Value = "  -0.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.52%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.79%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.97%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.109"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.31%  "

$ws.Range("D34").Value = "0.0₃0868"
$ws.Range("E34").Value = "  +7.11%  "

$ws.Range("E35").Value = "  +9.67%  "

$ws.Range("E36").Value = "  +2.79%  "

$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.37"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +13.51%  "

$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.14"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.27%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.80"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.42%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "449.19"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.50%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0374"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.67%  "

$ws.Range("D43").Value = "2.927.45"
$ws.Range("E43").Value = "  +5.00%  "

$ws.Range("E44").Value = "  +5.95%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.112"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.29%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "127.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.57%  "

$ws.Range("E48").Value = "  -0.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.112"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.81%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.68%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -10.06%  "
